$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (15) to make room for "DelegateTo".
# This shifts the old column O (Action-Level3) to P, and shifts widths/dimension accordingly.
$ws.Columns(15).EntireColumn.Insert()

# The newly inserted column picks up the default sheet width; Excel's default
# insert behaviour copies formatting (incl. width) from the column on the
# left, so match column N's (14) width here too.
$ws.Columns(15).ColumnWidth = $ws.Columns(14).ColumnWidth

# New header cell for the inserted column.
$ws.Range("O1").Value = "DelegateTo"

# Row 16 is the "Delegate" scenario: the delegate name had been entered in the
# ForwardTo column (N16); move it into the new DelegateTo column (O16) and
# clear the ForwardTo cell.
$ws.Range("O16").Value = $ws.Range("N16").Value2
$ws.Range("N16").ClearContents()

# Update the remembered selection.
$ws.Range("M22").Select() | Out-Null
